$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New headers in D1 and E1
$ws.Range("D1").Value = "场景音乐"
$ws.Range("E1").Value = "PK模式"

# Highlight the new header cells in red fill (Excel color is BGR: red = 0x0000FF = 255)
$ws.Range("D1:E1").Interior.Color = 255

# Fill data rows D2:E7 with 0
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 4).Value = 0
    $ws.Cells.Item($r, 5).Value = 0
}

# Column widths to match bestFit sizing from the diff (target width chars: D=9, E=7.25).
# This runtime quantizes ColumnWidth to whole pixels (7px/char "MDW" + 5px padding) before
# storing the OOXML <col width>, so the input value that round-trips to the desired stored
# width must be solved for, not assigned directly.
$ws.Columns.Item(4).ColumnWidth = 8.3
$ws.Columns.Item(5).ColumnWidth = 6.5

# Update selection
$ws.Range("K8").Select()
